$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range('D2').Value = '41.584.41'
$ws.Range('E2').Value = '  +0.13%  '
$ws.Range('D3').Value = '2.468.95'
$ws.Range('E4').Value = '  +0.27%  '
Set-TextValue $ws.Range('D5') '317.82'
$ws.Range('E5').Value = '  +1.48%  '
Set-TextValue $ws.Range('D6') '92.01'
$ws.Range('E6').Value = '  -0.44%  '
$ws.Range('E7').Value = '  +0.84%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('E9').Value = '  +0.57%  '
$ws.Range('B10').Value = 'Avalanche'
$ws.Range('C10').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range('D10') '32.88'
$ws.Range('E10').Value = '  +0.06%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws.Range('D11') '0.0852'
$ws.Range('E11').Value = '  +7.45%  '
$ws.Range('E12').Value = '  +0.42%  '
$ws.Range('D13').Value = '2.849.60'
$ws.Range('E13').Value = '  -0.46%  '
$ws.Range('E14').Value = '  -0.66%  '
Set-TextValue $ws.Range('D15') '15.52'
$ws.Range('E15').Value = '  -5.50%  '
$ws.Range('D16').Value = '2.474.68'
$ws.Range('E16').Value = '  -0.62%  '
Set-TextValue $ws.Range('D17') '0.790'
$ws.Range('E17').Value = '  +1.85%  '
$ws.Range('D18').Value = '41.538.73'
$ws.Range('E18').Value = '  -0.01%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range('D19') '6.44'
$ws.Range('E19').Value = '  -1.55%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0946'
$ws.Range('E20').Value = '  -0.09%  '
Set-TextValue $ws.Range('D21') '71.05'
$ws.Range('E21').Value = '  -1.76%  '
Set-TextValue $ws.Range('D22') '11.28'
$ws.Range('E22').Value = '  +0.57%  '
Set-TextValue $ws.Range('D23') '238.81'
$ws.Range('E23').Value = '  +0.86%  '
$ws.Range('E24').Value = '  +0.75%  '
Set-TextValue $ws.Range('D25') '1.92'
$ws.Range('E25').Value = '  +0.98%  '
Set-TextValue $ws.Range('D26') '1.00'
$ws.Range('E26').Value = '  -0.05%  '
Set-TextValue $ws.Range('D27') '24.58'
$ws.Range('E27').Value = '  -0.90%  '
Set-TextValue $ws.Range('D29') '9.83'
$ws.Range('E29').Value = '  +1.18%  '
Set-TextValue $ws.Range('D30') '36.07'
$ws.Range('E30').Value = '  +0.16%  '
Set-TextValue $ws.Range('D31') '161.20'
$ws.Range('E31').Value = '  +2.15%  '
$ws.Range('E32').Value = '  +0.46%  '
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('E34').Value = '  +0.52%  '
$ws.Range('E35').Value = '  +0.93%  '
Set-TextValue $ws.Range('D36') '17.23'
$ws.Range('E36').Value = '  -1.13%  '
$ws.Range('E37').Value = '  -0.46%  '
Set-TextValue $ws.Range('D38') '0.116'
$ws.Range('E38').Value = '  +1.32%  '
Set-TextValue $ws.Range('D39') '1.83'
$ws.Range('E39').Value = '  -0.31%  '
$ws.Range('E40').Value = '  -2.91%  '
Set-TextValue $ws.Range('D41') '3.98'
$ws.Range('E41').Value = '  -2.66%  '
Set-TextValue $ws.Range('D42') '2.45'
$ws.Range('E42').Value = '  +3.50%  '
$ws.Range('D43').Value = '1.982.65'
$ws.Range('E43').Value = '  +0.45%  '
$ws.Range('E44').Value = '  +0.24%  '
Set-TextValue $ws.Range('D45') '18.93'
$ws.Range('E45').Value = '  -1.64%  '
Set-TextValue $ws.Range('D46') '2.97'
$ws.Range('E46').Value = '  +0.39%  '
Set-TextValue $ws.Range('D47') '9.17'
$ws.Range('E47').Value = '  +2.47%  '
$ws.Range('D48').Value = '2.707.71'
$ws.Range('E48').Value = '  -0.46%  '
Set-TextValue $ws.Range('D49') '97.26'
$ws.Range('E49').Value = '  -0.69%  '
Set-TextValue $ws.Range('D50') '74.08'
$ws.Range('E50').Value = '  +2.17%  '
Set-TextValue $ws.Range('D51') '67.15'
$ws.Range('E51').Value = '  -1.71%  '
